$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.794.44'
$ws.Range('E2').Value = '  -2.45%  '
$ws.Range('D3').Value = '1.567.17'
$ws.Range('E3').Value = '  -0.14%  '
$ws.Range('E4').Value = '  +0.01%  '
$c = $ws.Range('D5')
$c.Formula = "=""206.49"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E5').Value = '  -1.13%  '
$ws.Range('E6').Value = '  -2.24%  '
$ws.Range('E7').Value = '  +0.02%  '
$c = $ws.Range('D8')
$c.Formula = "=""22.01"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E8').Value = '  -0.78%  '
$ws.Range('E9').Value = '  -0.80%  '
$ws.Range('E10').Value = '  -1.33%  '
$ws.Range('E11').Value = '  -0.34%  '
$ws.Range('D12').Value = '1.791.38'
$ws.Range('E12').Value = '  -0.05%  '
$ws.Range('D13').Value = '1.566.66'
$ws.Range('E13').Value = '  -0.70%  '
$ws.Range('E14').Value = '  -2.49%  '
$c = $ws.Range('D15')
$c.Formula = "=""0.514"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E15').Value = '  -0.93%  '
$ws.Range('D16').Value = '26.802.95'
$c = $ws.Range('D17')
$c.Formula = "=""61.39"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E17').Value = '  -3.75%  '
$c = $ws.Range('D18')
$c.Formula = "=""7.42"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E18').Value = '  +1.91%  '
$c = $ws.Range('D19')
$c.Formula = "=""215.14"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E19').Value = '  +0.33%  '
$ws.Range('D20').Value = '0.0₃0677'
$ws.Range('E20').Value = '  -2.08%  '
$ws.Range('E21').Value = '  +0.04%  '
$ws.Range('E22').Value = '  -0.06%  '
$c = $ws.Range('D23')
$c.Formula = "=""9.29"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E23').Value = '  -3.02%  '
$ws.Range('E24').Value = '  -0.77%  '
$c = $ws.Range('D25')
$c.Formula = "=""153.25"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E25').Value = '  +0.28%  '
$c = $ws.Range('D26')
$c.Formula = "=""6.74"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E26').Value = '  +0.39%  '
$ws.Range('E27').Value = '  -0.34%  '
$ws.Range('E28').Value = '  -0.02%  '
$ws.Range('E29').Value = '  -1.47%  '
$c = $ws.Range('D30')
$c.Formula = "=""0.0467"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E30').Value = '  -0.82%  '
$c = $ws.Range('D31')
$c.Formula = "=""1.11"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E31').Value = '  -3.37%  '
$ws.Range('E32').Value = '  -1.14%  '
$ws.Range('E33').Value = '  +1.45%  '
$ws.Range('E34').Value = '  -1.43%  '
$ws.Range('E35').Value = '  -1.13%  '
$ws.Range('E36').Value = '  -1.12%  '
$c = $ws.Range('D37')
$c.Formula = "=""0.935"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E37').Value = '  -1.72%  '
$ws.Range('E38').Value = '  -2.95%  '
$c = $ws.Range('D39')
$c.Formula = "=""0.529"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E39').Value = '  -2.79%  '
$c = $ws.Range('D40')
$c.Formula = "=""0.815"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E40').Value = '  -1.58%  '
$ws.Range('E41').Value = '  +0.05%  '
$c = $ws.Range('D42')
$c.Formula = "=""0.989"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E42').Value = '  +0.93%  '
$ws.Range('E43').Value = '  -0.09%  '
$ws.Range('E44').Value = '  +1.00%  '
$c = $ws.Range('D45')
$c.Formula = "=""5.32"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E45').Value = '  +0.89%  '
$c = $ws.Range('D46')
$c.Formula = "=""63.32"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E46').Value = '  -1.52%  '
$ws.Range('D47').Value = '1.703.70'
$ws.Range('E47').Value = '  +0.12%  '
$c = $ws.Range('D48')
$c.Formula = "=""85.99"""
$c.Copy()
$c.PasteSpecial(-4163)
$ws.Range('E48').Value = '  +0.76%  '
$ws.Range('D49').Value = '0.0₇0984'
$ws.Range('E49').Value = '  -1.42%  '
$ws.Range('E50').Value = '  -0.45%  '
$ws.Range('E51').Value = '  -0.93%  '
